$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $origStyle = $Range.Style
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '43.235.58'
Set-TextValue $ws.Range('E2') '  -1.17%  '
Set-TextValue $ws.Range('D3') '2.351.82'
Set-TextValue $ws.Range('E3') '  +4.42%  '
Set-TextValue $ws.Range('E4') '  -0.42%  '
Set-TextValue $ws.Range('D5') '233.78'
Set-TextValue $ws.Range('E5') '  +1.00%  '
Set-TextValue $ws.Range('E6') '  +0.27%  '
Set-TextValue $ws.Range('D7') '71.34'
Set-TextValue $ws.Range('E7') '  +12.81%  '
Set-TextValue $ws.Range('E8') '  -0.17%  '
Set-TextValue $ws.Range('D9') '0.485'
Set-TextValue $ws.Range('E9') '  +9.42%  '
Set-TextValue $ws.Range('D10') '0.0977'
Set-TextValue $ws.Range('E10') '  +1.34%  '
Set-TextValue $ws.Range('D11') '27.32'
Set-TextValue $ws.Range('E11') '  +3.06%  '
Set-TextValue $ws.Range('B12') 'WrappedliquidstakedEther2.0'
Set-TextValue $ws.Range('C12') 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range('D12') '2.706.38'
Set-TextValue $ws.Range('E12') '  +4.59%  '
Set-TextValue $ws.Range('B13') 'TRON'
Set-TextValue $ws.Range('C13') 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range('D13') '0.107'
Set-TextValue $ws.Range('E13') '  +1.23%  '
Set-TextValue $ws.Range('D14') '16.04'
Set-TextValue $ws.Range('E14') '  +3.20%  '
Set-TextValue $ws.Range('E15') '  +3.34%  '
Set-TextValue $ws.Range('E16') '  +3.37%  '
Set-TextValue $ws.Range('D17') '2.354.33'
Set-TextValue $ws.Range('E17') '  +4.40%  '
Set-TextValue $ws.Range('D18') '43.223.58'
Set-TextValue $ws.Range('E18') '  -0.93%  '
Set-TextValue $ws.Range('E19') '  +4.06%  '
Set-TextValue $ws.Range('D20') '6.33'
Set-TextValue $ws.Range('E20') '  +3.49%  '
Set-TextValue $ws.Range('D21') '74.40'
Set-TextValue $ws.Range('E21') '  +1.32%  '
Set-TextValue $ws.Range('D22') '249.97'
Set-TextValue $ws.Range('E22') '  +1.05%  '
Set-TextValue $ws.Range('E23') '  -0.03%  '
Set-TextValue $ws.Range('D24') '3.78'
Set-TextValue $ws.Range('E24') '  +1.73%  '
Set-TextValue $ws.Range('E25') '  +0.59%  '
Set-TextValue $ws.Range('B26') 'Toncoin'
Set-TextValue $ws.Range('C26') 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range('D26') '2.26'
Set-TextValue $ws.Range('E26') '  -1.55%  '
Set-TextValue $ws.Range('B27') 'Cosmos'
Set-TextValue $ws.Range('C27') 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range('D27') '10.02'
Set-TextValue $ws.Range('E27') '  +1.82%  '
Set-TextValue $ws.Range('D28') '22.34'
Set-TextValue $ws.Range('E28') '  +3.25%  '
Set-TextValue $ws.Range('D29') '172.75'
Set-TextValue $ws.Range('E29') '  +0.18%  '
Set-TextValue $ws.Range('E30') '  +7.21%  '
Set-TextValue $ws.Range('E31') '  -4.08%  '
Set-TextValue $ws.Range('E32') '  +0.61%  '
Set-TextValue $ws.Range('D33') '4.98'
Set-TextValue $ws.Range('E33') '  +1.81%  '
Set-TextValue $ws.Range('E34') '  +1.61%  '
Set-TextValue $ws.Range('D35') '5.06'
Set-TextValue $ws.Range('E35') '  +3.11%  '
Set-TextValue $ws.Range('E36') '  +2.14%  '
Set-TextValue $ws.Range('D37') '6.54'
Set-TextValue $ws.Range('D38') '2.43'
Set-TextValue $ws.Range('E38') '  +6.32%  '
Set-TextValue $ws.Range('E39') '  +0.50%  '
Set-TextValue $ws.Range('B40') 'InjectiveProtocol'
Set-TextValue $ws.Range('C40') 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D40') '18.80'
Set-TextValue $ws.Range('E40') '  +9.97%  '
Set-TextValue $ws.Range('B41') 'BinanceUSD'
Set-TextValue $ws.Range('C41') 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range('D41') '1.00'
Set-TextValue $ws.Range('E41') '  -0.15%  '
Set-TextValue $ws.Range('B42') 'FraxShare'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D42') '8.89'
Set-TextValue $ws.Range('E42') '  +3.64%  '
Set-TextValue $ws.Range('E43') '  +7.87%  '
Set-TextValue $ws.Range('D44') '4.48'
Set-TextValue $ws.Range('E44') '  -2.59%  '
Set-TextValue $ws.Range('D45') '99.02'
Set-TextValue $ws.Range('E45') '  +1.11%  '
Set-TextValue $ws.Range('E46') '  +2.10%  '
Set-TextValue $ws.Range('D47') '0.0958'
Set-TextValue $ws.Range('E47') '  +1.97%  '
Set-TextValue $ws.Range('D48') '1.440.41'
Set-TextValue $ws.Range('E48') '  -0.67%  '
Set-TextValue $ws.Range('D49') '2.581.46'
Set-TextValue $ws.Range('E49') '  +4.77%  '
Set-TextValue $ws.Range('E50') '  +0.28%  '
Set-TextValue $ws.Range('E51') '  -2.50%  '
